$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1298.7142
$ws.Range("I6").Value = 220.4
$ws.Range("J6").Value = 3994.5
$ws.Range("K6").Value = 661.2
$ws.Range("L6").Value = 11983.5
$ws.Range("M6").Value = -549.2
$ws.Range("N6").Value = -12207.5
$ws.Range("H129").Value = 37037756
$ws.Range("I129").Value = 100000330
$ws.Range("J129").Value = 952.8823
$ws.Range("K129").Value = 300000990
$ws.Range("L129").Value = 2858.6469
$ws.Range("M129").Value = -299995990
$ws.Range("N129").Value = -12858.6469
$ws.Range("H132").Value = 2969
$ws.Range("I132").Value = 3076.8333
$ws.Range("K132").Value = 9230.499899999999
$ws.Range("M132").Value = -6700.499899999999
$ws.Range("H135").Value = 580.0769
$ws.Range("I135").Value = 241.36363
$ws.Range("K135").Value = 2172.27267
$ws.Range("M135").Value = 362.7273300000002
$ws.Range("H138").Value = 2695.6875
$ws.Range("I138").Value = 1192.7894
$ws.Range("J138").Value = 4892.231
$ws.Range("K138").Value = 3578.3682
$ws.Range("L138").Value = 14676.693
$ws.Range("M138").Value = 1561.6318
$ws.Range("N138").Value = -24956.693
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4575.9805
$ws.Range("I32").Value = 2586.3333
$ws.Range("K32").Value = 2586.3333
$ws.Range("M32").Value = -2299.3333
$ws.Range("H61").Value = 20836854
$ws.Range("I61").Value = 26318274
$ws.Range("J61").Value = 7460.9
$ws.Range("K61").Value = 26318274
$ws.Range("L61").Value = 7460.9
$ws.Range("M61").Value = -26318062
$ws.Range("N61").Value = -7884.9
$ws.Range("H136").Value = 20836854
$ws.Range("I136").Value = 26318274
$ws.Range("J136").Value = 7460.9
$ws.Range("K136").Value = 78954822
$ws.Range("L136").Value = 22382.7
$ws.Range("M136").Value = -78952272
$ws.Range("N136").Value = -27482.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H105").Value = 2560.7
$ws.Range("I105").Value = 2623
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2623
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -876
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 2012.1702
$ws.Range("I134").Value = 1330.1052
$ws.Range("K134").Value = 3990.3156
$ws.Range("M134").Value = -1455.3156
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10300.954
$ws.Range("I31").Value = 1411.6538
$ws.Range("J31").Value = 23141.055
$ws.Range("K31").Value = 1411.6538
$ws.Range("L31").Value = 23141.055
$ws.Range("M31").Value = -1116.6538
$ws.Range("N31").Value = -23731.055
$ws.Range("H34").Value = 10300.954
$ws.Range("I34").Value = 1411.6538
$ws.Range("J34").Value = 23141.055
$ws.Range("K34").Value = 1411.6538
$ws.Range("L34").Value = 23141.055
$ws.Range("M34").Value = -1209.6538
$ws.Range("N34").Value = -23545.055
$ws.Range("H134").Value = 4047
$ws.Range("I134").Value = 3653.2942
$ws.Range("J134").Value = 5385.6
$ws.Range("K134").Value = 10959.8826
$ws.Range("L134").Value = 16156.8
$ws.Range("M134").Value = -8424.882599999999
$ws.Range("N134").Value = -21226.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2012.1177
$ws.Range("I5").Value = 968
$ws.Range("J5").Value = 2333.3845
$ws.Range("K5").Value = 2904
$ws.Range("L5").Value = 7000.1535
$ws.Range("M5").Value = -2792
$ws.Range("N5").Value = -7224.1535
$ws.Range("H122").Value = 713.0476
$ws.Range("J122").Value = 1083
$ws.Range("L122").Value = 9747
$ws.Range("N122").Value = -14647
$ws.Range("H131").Value = 1403.326
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1403.326
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4209.978
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -14289.978
$ws.Range("H135").Value = 2012.1177
$ws.Range("I135").Value = 968
$ws.Range("J135").Value = 2333.3845
$ws.Range("K135").Value = 8712
$ws.Range("L135").Value = 21000.4605
$ws.Range("M135").Value = -6177
$ws.Range("N135").Value = -26070.4605
$ws.Range("H138").Value = 1427.8948
$ws.Range("I138").Value = 913
$ws.Range("K138").Value = 2739
$ws.Range("M138").Value = 2401
$ws.Range("H141").Value = 2893.125
$ws.Range("I141").Value = 548.3333
$ws.Range("J141").Value = 4300
$ws.Range("K141").Value = 1644.9999
$ws.Range("L141").Value = 12900
$ws.Range("M141").Value = 3535.0001
$ws.Range("N141").Value = -23260
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H107").Value = 1709.3846
$ws.Range("I107").Value = 2132.9
$ws.Range("J107").Value = 297.66666
$ws.Range("K107").Value = 2132.9
$ws.Range("L107").Value = 297.66666
$ws.Range("M107").Value = -212.9000000000001
$ws.Range("N107").Value = -4137.66666
$ws.Range("H122").Value = 8283.4
$ws.Range("I122").Value = 1634.0769
$ws.Range("J122").Value = 51504
$ws.Range("K122").Value = 4902.2307
$ws.Range("L122").Value = 154512
$ws.Range("M122").Value = -2452.2307
$ws.Range("N122").Value = -159412
$ws.Range("H132").Value = 3170.3704
$ws.Range("I132").Value = 2049
$ws.Range("J132").Value = 4067.4666
$ws.Range("K132").Value = 6147
$ws.Range("L132").Value = 12202.3998
$ws.Range("M132").Value = -3617
$ws.Range("N132").Value = -17262.3998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 126512.625
$ws.Range("I7").Value = 167617
$ws.Range("J7").Value = 3199.5
$ws.Range("K7").Value = 167617
$ws.Range("L7").Value = 3199.5
$ws.Range("M7").Value = -167505
$ws.Range("N7").Value = -3423.5
$ws.Range("H46").Value = 23858242
$ws.Range("I46").Value = 55667330
$ws.Range("J46").Value = 1427.25
$ws.Range("K46").Value = 55667330
$ws.Range("L46").Value = 1427.25
$ws.Range("M46").Value = -55667142
$ws.Range("N46").Value = -1803.25
$ws.Range("H61").Value = 1923.3334
$ws.Range("I61").Value = 1544
$ws.Range("J61").Value = 2397.5
$ws.Range("K61").Value = 1544
$ws.Range("L61").Value = 2397.5
$ws.Range("M61").Value = -1342
$ws.Range("N61").Value = -2801.5
$ws.Range("H93").Value = 15196.143
$ws.Range("J93").Value = 780
$ws.Range("L93").Value = 780
$ws.Range("N93").Value = -3276
$ws.Range("H113").Value = 1923.3334
$ws.Range("I113").Value = 1544
$ws.Range("J113").Value = 2397.5
$ws.Range("K113").Value = 1544
$ws.Range("L113").Value = 2397.5
$ws.Range("M113").Value = 626
$ws.Range("N113").Value = -6737.5
$ws.Range("H126").Value = 126512.625
$ws.Range("I126").Value = 167617
$ws.Range("J126").Value = 3199.5
$ws.Range("K126").Value = 502851
$ws.Range("L126").Value = 9598.5
$ws.Range("M126").Value = -500381
$ws.Range("N126").Value = -14538.5
$ws.Range("H136").Value = 3332.3044
$ws.Range("I136").Value = 2096.3845
$ws.Range("K136").Value = 6289.1535
$ws.Range("M136").Value = -3739.1535
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 211388
$ws.Range("J46").Value = 211388
$ws.Range("L46").Value = 211388
$ws.Range("N46").Value = -211850
$ws.Range("H107").Value = 598.5
$ws.Range("I107").Value = 517.125
$ws.Range("J107").Value = 761.25
$ws.Range("K107").Value = 1551.375
$ws.Range("L107").Value = 2283.75
$ws.Range("M107").Value = 368.625
$ws.Range("N107").Value = -6123.75
$ws.Range("H132").Value = 2130.439
$ws.Range("I132").Value = 1561.6666
$ws.Range("J132").Value = 2727.65
$ws.Range("K132").Value = 4684.9998
$ws.Range("L132").Value = 8182.950000000001
$ws.Range("M132").Value = -2154.9998
$ws.Range("N132").Value = -13242.95
$ws.Range("H133").Value = 24833.334
$ws.Range("J133").Value = 24833.334
$ws.Range("L133").Value = 24833.334
$ws.Range("N133").Value = -34953.334
$ws.Range("H134").Value = 211388
$ws.Range("J134").Value = 211388
$ws.Range("L134").Value = 634164
$ws.Range("N134").Value = -639234
$ws.Range("H136").Value = 4412.939
$ws.Range("I136").Value = 776.125
$ws.Range("J136").Value = 5122.561
$ws.Range("K136").Value = 2328.375
$ws.Range("L136").Value = 15367.683
$ws.Range("M136").Value = 221.625
$ws.Range("N136").Value = -20467.683

Write-Host "edit complete"